# This script applies a weekly price-sheet update to the "Cebolla" (onion)
# price log. Two new price observations are inserted at the top of the
# existing data block (rows 921-922), pushing all the subsequent rows
# down by two (the sheet's dimension grows from A1:R964 to A1:R966).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows right before the existing row 921, shifting
# everything from row 921 down onward down by two rows.
$ws.Range("A921:A922").EntireRow.Insert()

# --- New row 921 ---
$ws.Range("A921").Value = 7
$ws.Range("B921").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C921").Value = "Ñuble"
$ws.Range("D921").Value = 45041
$ws.Range("E921").Value = 16
$ws.Range("F921").Value = 100112004
$ws.Range("G921").Value = "Cebolla"
$ws.Range("H921").Value = "Sin especificar"
$ws.Range("I921").Value = "1a (cosecha)"
$ws.Range("J921").Value = 120
$ws.Range("K921").Value = 8000
$ws.Range("L921").Value = 8500
$ws.Range("M921").Value = 8250
$ws.Range("N921").Value = "$/malla 18 kilos"
$ws.Range("O921").Value = "Región de O'Higgins"
$ws.Range("P921").Value = 458
$ws.Range("Q921").Value = 18
$ws.Range("R921").Value = "Hortaliza"

# --- New row 922 ---
$ws.Range("A922").Value = 7
$ws.Range("B922").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C922").Value = "Ñuble"
$ws.Range("D922").Value = 45041
$ws.Range("E922").Value = 16
$ws.Range("F922").Value = 100112004
$ws.Range("G922").Value = "Cebolla"
$ws.Range("H922").Value = "Sin especificar"
$ws.Range("I922").Value = "1a (cosecha)"
$ws.Range("J922").Value = 300
$ws.Range("K922").Value = 11000
$ws.Range("L922").Value = 11000
$ws.Range("M922").Value = 11000
$ws.Range("N922").Value = "$/malla 25 kilos"
$ws.Range("O922").Value = "Región del Maule"
$ws.Range("P922").Value = 440
$ws.Range("Q922").Value = 25
$ws.Range("R922").Value = "Hortaliza"
